# Generate Report for Handoff
# Update status text from "In Translation" to "Ready for handoff" and
# bump the handoff timestamps on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and generate date (G2)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-18 02:51:43"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-18 02:51:39"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-18 02:51:43"

# Widen the status/date columns so their widths match the longer text
# (mirrors Excel's auto-fit after the longer "Ready for handoff" string)
$overview.Range("E1").ColumnWidth = 16.4
$overview.Range("F1").ColumnWidth = 16.4
$zhcn.Range("C1").ColumnWidth = 16.4
$dede.Range("C1").ColumnWidth = 16.4
